$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.445.25"
$ws.Range("E2").Value = "  -1.22%  "
$ws.Range("D3").Value = "'1.646.86"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'1.001"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").Value = "'298.73"
$ws.Range("E6").Value = "  -1.65%  "
$ws.Range("D7").Value = "'0.3785"
$ws.Range("E7").Value = "  -1.03%  "
$ws.Range("D8").Value = "'0.3556"
$ws.Range("E8").Value = "  -1.44%  "
$ws.Range("D9").Value = "'49.85"
$ws.Range("E9").Value = "  -2.78%  "
$ws.Range("D10").Value = "'0.08091"
$ws.Range("E10").Value = "  -1.75%  "
$ws.Range("D11").Value = "'1.220"
$ws.Range("E11").Value = "  -2.40%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").Value = "'22.05"
$ws.Range("E13").Value = "  -2.56%  "
$ws.Range("D14").Value = "'6.394"
$ws.Range("E14").Value = "  -2.29%  "
$ws.Range("D15").Value = "'7.345"
$ws.Range("E15").Value = "  -0.86%  "
$ws.Range("D16").Value = "'0.00001197"
$ws.Range("E16").Value = "  -3.08%  "
$ws.Range("D17").Value = "'1.646.52"
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("D18").Value = "'97.47"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").Value = "'0.06951"
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("D20").Value = "'6.759"
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("D21").Value = "'17.34"
$ws.Range("E21").Value = "  -2.01%  "
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").Value = "'12.41"
$ws.Range("E23").Value = "  -1.75%  "
$ws.Range("D24").Value = "'23.451.15"
$ws.Range("E24").Value = "  -1.21%  "
$ws.Range("D25").Value = "'2.493"
$ws.Range("E25").Value = "  -1.17%  "
$ws.Range("D26").Value = "'2.901"
$ws.Range("E26").Value = "  -6.07%  "
$ws.Range("D27").Value = "'20.92"
$ws.Range("E27").Value = "  -1.89%  "
$ws.Range("D28").Value = "'152.62"
$ws.Range("E28").Value = "  +0.74%  "
$ws.Range("D29").Value = "'5.214"
$ws.Range("E29").Value = "  -1.18%  "
$ws.Range("D30").Value = "'132.70"
$ws.Range("E30").Value = "  -1.65%  "
$ws.Range("D31").Value = "'1.832.27"
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("D32").Value = "'6.911"
$ws.Range("E32").Value = "  +0.63%  "
$ws.Range("B33").Value = "WEMIXTOKEN"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").Value = "'2.118"
$ws.Range("E33").Value = "  +0.54%  "
$ws.Range("B34").Value = "FraxShare"
$ws.Range("C34").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D34").Value = "'11.66"
$ws.Range("E34").Value = "  -1.45%  "
$ws.Range("D35").Value = "'1.002"
$ws.Range("E35").Value = "  -8.34%  "
$ws.Range("D36").Value = "'0.02715"
$ws.Range("E36").Value = "  -4.55%  "
$ws.Range("D37").Value = "'0.08729"
$ws.Range("E37").Value = "  -1.15%  "
$ws.Range("D38").Value = "'0.2424"
$ws.Range("E38").Value = "  -3.73%  "
$ws.Range("D39").Value = "'5.919"
$ws.Range("E39").Value = "  -3.02%  "
$ws.Range("D40").Value = "'12.98"
$ws.Range("E40").Value = "  +0.87%  "
$ws.Range("D41").Value = "'0.06776"
$ws.Range("E41").Value = "  -3.94%  "
$ws.Range("D42").Value = "'0.6873"
$ws.Range("E42").Value = "  -2.90%  "
$ws.Range("D43").Value = "'1.303"
$ws.Range("E43").Value = "  -2.65%  "
$ws.Range("D44").Value = "'15.53"
$ws.Range("E44").Value = "  -3.10%  "
$ws.Range("E45").Value = "  +0.19%  "
$ws.Range("D46").Value = "'0.6360"
$ws.Range("E46").Value = "  -2.92%  "
$ws.Range("D47").Value = "'2.248"
$ws.Range("E47").Value = "  -3.94%  "
$ws.Range("D48").Value = "'3.912"
$ws.Range("E48").Value = "  -1.50%  "
$ws.Range("D49").Value = "'0.07718"
$ws.Range("E49").Value = "  -3.35%  "
$ws.Range("D50").Value = "'127.09"
$ws.Range("E50").Value = "  -1.13%  "
$ws.Range("D51").Value = "'1.149"
$ws.Range("E51").Value = "  -3.97%  "
